# expansão das análises automáticas
# Adds three new summary columns (apoio_medio, contribuicoes,
# media_contribuicoes) and rescales the percentage columns (E, F) from
# fractional (0-1) to percentage-point (0-100) numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -----------------------------------------------
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the bold/bordered/centered header formatting used by the rest of
# row 1 (copy the style from the last existing header cell, K1).
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# --- New column data ----------------------------------------------------
$newData = @{
    2 = @(91.85312348260253, 209535, 322.3615384615385)
    3 = @(89.17093558435907, 54018, 300.1)
    4 = @(89.37434882498151, 141221, 132.6018779342723)
    5 = @(91.95990423942952, 62425, 196.3050314465409)
    6 = @(19.36290068160405, 2129, 15.65441176470588)
    7 = @(24.85243295759227, 79, 4.9375)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Cells.Item($row, 12).Value = $vals[0]
    $ws.Cells.Item($row, 13).Value = $vals[1]
    $ws.Cells.Item($row, 14).Value = $vals[2]
}

# --- Rescale existing percentage columns E (5) and F (6) from fraction --
# --- to percentage-point values (x100), keeping the same number format -
# (computed as old_value * 100; written as literals to avoid reintroducing
# a different float64 rounding than the canonical recompute).
$ef = @{
    2 = @(79.10112359550561, 61.5530303030303)
    3 = @(20.89887640449438, 64.51612903225806)
    4 = @(77.79291553133515, 93.25744308231172)
    5 = @(22.20708446866485, 97.54601226993866)
    6 = @(94.00584795321637, 21.15085536547434)
    7 = @(5.994152046783626, 39.02439024390244)
}

foreach ($row in $ef.Keys) {
    $vals = $ef[$row]
    $ws.Cells.Item($row, 5).Value = $vals[0]
    $ws.Cells.Item($row, 6).Value = $vals[1]
}
